{"js": "// Replace the date line and each \"NNN\u00d7N=\" equation cell with its new value.\n// The document has one occurrence of each \"from\" string, so body.search()\n// with matchCase (and exact text) safely targets a single run each time.\nconst replacements = [\n  [\"2024-02-14 Wednesday\", \"2024-02-15 Thursday\"],\n  [\"258\u00d76=\", \"490\u00d75=\"],\n  [\"247\u00d72=\", \"512\u00d75=\"],\n  [\"363\u00d78=\", \"155\u00d76=\"],\n  [\"974\u00d79=\", \"864\u00d75=\"],\n  [\"693\u00d78=\", \"395\u00d74=\"],\n  [\"587\u00d73=\", \"422\u00d76=\"],\n  [\"977\u00d75=\", \"520\u00d76=\"],\n  [\"599\u00d78=\", \"432\u00d76=\"],\n  [\"216\u00d79=\", \"895\u00d79=\"],\n  [\"902\u00d76=\", \"825\u00d76=\"],\n  [\"824\u00d77=\", \"646\u00d75=\"],\n  [\"460\u00d79=\", \"151\u00d75=\"],\n  [\"833\u00d73=\", \"888\u00d72=\"],\n  [\"311\u00d73=\", \"938\u00d74=\"],\n  [\"803\u00d78=\", \"839\u00d77=\"],\n  [\"912\u00d79=\", \"189\u00d78=\"],\n  [\"223\u00d79=\", \"160\u00d77=\"],\n  [\"987\u00d73=\", \"650\u00d79=\"],\n  [\"318\u00d73=\", \"683\u00d77=\"],\n  [\"799\u00d74=\", \"254\u00d73=\"],\n  [\"293\u00d75=\", \"114\u00d75=\"],\n  [\"690\u00d73=\", \"348\u00d74=\"],\n  [\"748\u00d76=\", \"587\u00d74=\"],\n  [\"956\u00d76=\", \"810\u00d78=\"],\n  [\"350\u00d73=\", \"752\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${from}`);\n  }\n  // Replace only the first occurrence; each \"from\" string is unique in the\n  // document, so this is always the correct (and only) match.\n  results.items[0].insertText(to, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=\" equation cell with its new value.\n# Every \"from\" string occurs exactly once in the document, so a scoped\n# Find/Replace (Replace = wdReplaceOne, no wraparound needed since we start\n# from the top each time) always rewrites the single matching run in place,\n# preserving that run's formatting (font/size).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-14 Wednesday\", \"2024-02-15 Thursday\"),\n    @(\"258\u00d76=\", \"490\u00d75=\"),\n    @(\"247\u00d72=\", \"512\u00d75=\"),\n    @(\"363\u00d78=\", \"155\u00d76=\"),\n    @(\"974\u00d79=\", \"864\u00d75=\"),\n    @(\"693\u00d78=\", \"395\u00d74=\"),\n    @(\"587\u00d73=\", \"422\u00d76=\"),\n    @(\"977\u00d75=\", \"520\u00d76=\"),\n    @(\"599\u00d78=\", \"432\u00d76=\"),\n    @(\"216\u00d79=\", \"895\u00d79=\"),\n    @(\"902\u00d76=\", \"825\u00d76=\"),\n    @(\"824\u00d77=\", \"646\u00d75=\"),\n    @(\"460\u00d79=\", \"151\u00d75=\"),\n    @(\"833\u00d73=\", \"888\u00d72=\"),\n    @(\"311\u00d73=\", \"938\u00d74=\"),\n    @(\"803\u00d78=\", \"839\u00d77=\"),\n    @(\"912\u00d79=\", \"189\u00d78=\"),\n    @(\"223\u00d79=\", \"160\u00d77=\"),\n    @(\"987\u00d73=\", \"650\u00d79=\"),\n    @(\"318\u00d73=\", \"683\u00d77=\"),\n    @(\"799\u00d74=\", \"254\u00d73=\"),\n    @(\"293\u00d75=\", \"114\u00d75=\"),\n    @(\"690\u00d73=\", \"348\u00d74=\"),\n    @(\"748\u00d76=\", \"587\u00d74=\"),\n    @(\"956\u00d76=\", \"810\u00d78=\"),\n    @(\"350\u00d73=\", \"752\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $to\n\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $found = $find.Execute($from, $false, $false, $false, $false, $false, $true, 1, $false, $to, 1)\n    if (-not $found) {\n        throw \"Text not found: $from\"\n    }\n}\n"}
